# Klipper US workbook update: add "US 7" sprint sheet content, fix a few
# task statuses/owners on existing sheets, and re-point the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "US7" worksheet right before "Issues" and populate it
#    with the Server Environment setup / deployment task list.
# ---------------------------------------------------------------------
$issues = $wb.Worksheets.Item("Issues")
$us7 = $wb.Worksheets.Add($issues, $null)
$us7.Name = "US7"

$us7.Range("B3").Value = "US 7"
$us7.Range("C3").Value = "Deployment"

$us7.Range("B5").Value = "Task No."
$us7.Range("C5").Value = "Tasks"
$us7.Range("D5").Value = "Remaining Time"
$us7.Range("E5").Value = "Assigned To"
$us7.Range("F5").Value = "Status"

$us7.Range("B6").Value = 1
$us7.Range("C6").Value = "Server Environment setup"
$us7.Range("E6").Value = "Sanket"
$us7.Range("F6").Value = "Complete"

$us7.Range("B7").Value = 2
$r = $us7.Range("C7")
$r.Value = "Server Environment setup-ProductionMachine"
$r.Characters(1, 26).Font.Italic = $false
$r.Characters(27, 17).Font.Italic = $true
$us7.Range("E7").Value = "Sanket"
$us7.Range("F7").Value = "Complete"

$us7.Range("B8").Value = 3
$r = $us7.Range("C8")
$r.Value = "server accessible within network-ProductionMachine"
$r.Characters(1, 33).Font.Italic = $false
$r.Characters(34, 17).Font.Italic = $true
$us7.Range("E8").Value = "IT"
$us7.Range("F8").Value = "InProgress"

$us7.Range("B9").Value = 4
$us7.Range("C9").Value = "Change in script(Read data for last 2 days only)"
$us7.Range("E9").Value = "Sanket"
$us7.Range("F9").Value = "Complete"

$us7.Range("B10").Value = 5
$us7.Range("C10").Value = "Deploy data-sync-schedular script on server"
$us7.Range("E10").Value = "Sanket"
$us7.Range("F10").Value = "Complete"

$us7.Range("B11").Value = 6
$us7.Range("C11").Value = "Change in script(sync extra fields from Access to Mongo)"
$us7.Range("E11").Value = "Sanket"
$us7.Range("F11").Value = "Complete"

$us7.Range("B12").Value = 7
$us7.Range("C12").Value = "Make script available on github"
$us7.Range("E12").Value = "Sanket"
$us7.Range("F12").Value = "Complete"

$us7.Range("B13").Value = 8
$us7.Range("C13").Value = "Add logging mechanism to script"
$us7.Range("E13").Value = "Sanket"
$us7.Range("F13").Value = "ToDo"

$us7.Range("B14").Value = 9
$us7.Range("C14").Value = "Data comparison and monitoring for a week"
$us7.Range("E14").Value = "Sanket"
$us7.Range("F14").Value = "InProgress"

$us7.Range("B3:F19").Select()
$us7.Range("E27").Select()

# ---------------------------------------------------------------------
# 2. US 1 - add row 11 "Add test cases" task.
# ---------------------------------------------------------------------
$us1 = $wb.Worksheets.Item("US 1")
$us1.Range("B11").Value = 7
$us1.Range("C11").Value = "Add test cases"
$us1.Range("D11").Value = 0
$us1.Range("E11").Value = "Sidhdesh"
$us1.Range("F11").Value = "Completed"
$us1.Range("F11").Select()

# ---------------------------------------------------------------------
# 3. US 2 - clear remaining-time for finished tasks, add row 9 "Test cases".
# ---------------------------------------------------------------------
$us2 = $wb.Worksheets.Item("US 2")
$us2.Range("D5").Value = 0
$us2.Range("D7").Value = 0
$us2.Range("D8").Value = 0

$us2.Range("B9").Value = 5
$us2.Range("C9").Value = "Test cases"
$us2.Range("D9").Value = 0
$us2.Range("E9").Value = "Sidhdesh"
$us2.Range("F9").Value = "Completed"

# ---------------------------------------------------------------------
# 4. US 3 - rename a task/owner, fix a typo, move the active-tab flag away.
# ---------------------------------------------------------------------
$us3 = $wb.Worksheets.Item("US 3")
$us3.Range("C8").Value = "Integration"
$us3.Range("E8").Value = "SIdhdesh"
$us3.Range("C10").Value = "Visibilit of reportee tab for lead and admin only"
$us3.Range("D9").Select()

# ---------------------------------------------------------------------
# 5. US 5 - clear the "Assigned To" value for the in-progress task.
# ---------------------------------------------------------------------
$us5 = $wb.Worksheets.Item("US 5")
$us5.Range("E7").ClearContents()
$us5.Range("E7").Select()

# ---------------------------------------------------------------------
# 6. Misc selection bookkeeping on the other sheets, to mirror what the
#    author had on-screen when the workbook was last saved.
# ---------------------------------------------------------------------
$usList = $wb.Worksheets.Item("US List")
$usList.Range("C14").Select()

$us4 = $wb.Worksheets.Item("US 4")
$us4.Range("E8").Select()

$us6 = $wb.Worksheets.Item("US 6")
$us6.Range("C32").Select()

$us7.Activate()
$us7.Range("E27").Select()

Write-Host "applied edits"
